$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I24").Value = "yields klaar"
$ws.Range("J24").Value = "1 uur"
$ws.Range("K24").NumberFormat = "d-mmm"
$ws.Range("K24").Value = 43020
$ws.Range("L24").Value = "was heel makkelijk"

$ws.Range("I25").Value = "begonnen met development page"
$ws.Range("J25").Value = "/"
$ws.Range("K25").Value = "/"
$ws.Range("L25").Value = "/"

$ws.Range("K28").Select()
